# Auto-generated edit script: updates currentAveragePrice/Leve price/profit
# columns (H-N) across multiple crafting-job worksheets to reflect refreshed
# market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 320.58334
$ws.Range("I38").Value = 49.666668
$ws.Range("K38").Value = 149.000004
$ws.Range("M38").Value = 222.999996
$ws.Range("H62").Value = 998.3333
$ws.Range("I62").Value = 998.3333
$ws.Range("K62").Value = 998.3333
$ws.Range("M62").Value = -374.3333
$ws.Range("H65").Value = 998.3333
$ws.Range("I65").Value = 998.3333
$ws.Range("K65").Value = 4991.6665
$ws.Range("M65").Value = -1871.6665
$ws.Range("H74").Value = 5999.5
$ws.Range("J74").Value = 5499
$ws.Range("L74").Value = 5499
$ws.Range("N74").Value = -7371
$ws.Range("H77").Value = 5999.5
$ws.Range("J77").Value = 5499
$ws.Range("L77").Value = 27495
$ws.Range("N77").Value = -36855
$ws.Range("H97").Value = 2455.889
$ws.Range("J97").Value = 2455.889
$ws.Range("L97").Value = 7367.667
$ws.Range("N97").Value = -8359.667000000001
$ws.Range("H99").Value = 542.8
$ws.Range("I99").Value = 616
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 1848
$ws.Range("L99").Value = 750
$ws.Range("M99").Value = -350
$ws.Range("N99").Value = -3746
$ws.Range("H106").Value = 4197.25
$ws.Range("I106").Value = 4197.25
$ws.Range("K106").Value = 4197.25
$ws.Range("M106").Value = -3566.25
$ws.Range("H137").Value = 1200
$ws.Range("J137").Value = 1200
$ws.Range("L137").Value = 3600
$ws.Range("N137").Value = -8700
$ws.Range("H141").Value = 1599
$ws.Range("I141").Value = 1599
$ws.Range("K141").Value = 4797
$ws.Range("M141").Value = 383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 822.8
$ws.Range("I110").Value = 822.8
$ws.Range("K110").Value = 822.8
$ws.Range("M110").Value = 1222.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1002.2
$ws.Range("J64").Value = 1004.3333
$ws.Range("L64").Value = 1004.3333
$ws.Range("N64").Value = -1454.3333
$ws.Range("H67").Value = 1002.2
$ws.Range("J67").Value = 1004.3333
$ws.Range("L67").Value = 1004.3333
$ws.Range("N67").Value = -2564.3333
$ws.Range("H86").Value = 1353.5
$ws.Range("I86").Value = 1955.4286
$ws.Range("J86").Value = 510.8
$ws.Range("K86").Value = 1955.4286
$ws.Range("L86").Value = 510.8
$ws.Range("M86").Value = -832.4286
$ws.Range("N86").Value = -2756.8
$ws.Range("H89").Value = 1353.5
$ws.Range("I89").Value = 1955.4286
$ws.Range("J89").Value = 510.8
$ws.Range("K89").Value = 9777.143
$ws.Range("L89").Value = 2554
$ws.Range("M89").Value = -4161.143
$ws.Range("N89").Value = -13786
$ws.Range("H94").Value = 1616.7142
$ws.Range("I94").Value = 1761.1666
$ws.Range("K94").Value = 1761.1666
$ws.Range("M94").Value = -1310.1666
$ws.Range("H99").Value = 3184.8
$ws.Range("I99").Value = 3094.3333
$ws.Range("K99").Value = 3094.3333
$ws.Range("M99").Value = -1596.3333
$ws.Range("H107").Value = 625.8
$ws.Range("I107").Value = 623.2414
$ws.Range("K107").Value = 623.2414
$ws.Range("M107").Value = 1296.7586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2448.3333
$ws.Range("I107").Value = 1890.2727
$ws.Range("J107").Value = 3325.2856
$ws.Range("K107").Value = 1890.2727
$ws.Range("L107").Value = 3325.2856
$ws.Range("M107").Value = 29.72730000000001
$ws.Range("N107").Value = -7165.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 70899.734
$ws.Range("I4").Value = 349.6
$ws.Range("J4").Value = 212000
$ws.Range("K4").Value = 1048.8
$ws.Range("L4").Value = 636000
$ws.Range("M4").Value = -936.8000000000002
$ws.Range("N4").Value = -636224
$ws.Range("H69").Value = 5665
$ws.Range("H72").Value = 5665
$ws.Range("H74").Value = 7000
$ws.Range("J74").Value = 7000
$ws.Range("L74").Value = 21000
$ws.Range("N74").Value = -23122
$ws.Range("H77").Value = 7000
$ws.Range("J77").Value = 7000
$ws.Range("L77").Value = 63000
$ws.Range("N77").Value = -73608
$ws.Range("H130").Value = 14750
$ws.Range("J130").Value = 15000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H131").Value = 950.7692
$ws.Range("I131").Value = 680
$ws.Range("J131").Value = 1032
$ws.Range("K131").Value = 2040
$ws.Range("L131").Value = 3096
$ws.Range("M131").Value = 3000
$ws.Range("N131").Value = -13176

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4359.6
$ws.Range("I70").Value = 4359.6
$ws.Range("K70").Value = 4359.6
$ws.Range("M70").Value = -4089.6
$ws.Range("H73").Value = 4359.6
$ws.Range("I73").Value = 4359.6
$ws.Range("K73").Value = 4359.6
$ws.Range("M73").Value = -3423.6
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3391.4
$ws.Range("I126").Value = 2834.6667
$ws.Range("K126").Value = 8504.000100000001
$ws.Range("M126").Value = -6034.000100000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 362
$ws.Range("I55").Value = 150.5
$ws.Range("J55").Value = 404.3
$ws.Range("K55").Value = 150.5
$ws.Range("L55").Value = 404.3
$ws.Range("M55").Value = 22.5
$ws.Range("N55").Value = -750.3
$ws.Range("H99").Value = 15859.8
$ws.Range("I99").Value = 15859.8
$ws.Range("K99").Value = 15859.8
$ws.Range("M99").Value = -12864.8
$ws.Range("H122").Value = 4236.375
$ws.Range("I122").Value = 4236.375
$ws.Range("K122").Value = 12709.125
$ws.Range("M122").Value = -10259.125
$ws.Range("H132").Value = 4674.8
$ws.Range("I132").Value = 4218.625
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 12655.875
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -10125.875
$ws.Range("N132").Value = -24558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 24950
$ws.Range("I70").Value = 24900
$ws.Range("K70").Value = 24900
$ws.Range("M70").Value = -24585
$ws.Range("H73").Value = 24950
$ws.Range("I73").Value = 24900
$ws.Range("K73").Value = 24900
$ws.Range("M73").Value = -23808
$ws.Range("H126").Value = 3250
$ws.Range("I126").Value = 3250
$ws.Range("K126").Value = 9750
$ws.Range("M126").Value = -7280
